# Issue #280 - Modify electricity sector CCS tracking
# Sheet "CPPbES" (CPP CO2 Capture Potential by Electricity Source)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPPbES")

# Electricity sources that no longer get a blanket 100% CO2-capture share
# (set to 0 - CCS not assumed available for these source types anymore)
$ws.Range("B2").Value  = 0   # hard coal
$ws.Range("B3").Value  = 0   # natural gas steam turbine
$ws.Range("B4").Value  = 0   # natural gas combined cycle
$ws.Range("B10").Value = 0   # biomass
$ws.Range("B12").Value = 0   # petroleum
$ws.Range("B13").Value = 0   # natural gas peaker
$ws.Range("B14").Value = 0   # lignite
$ws.Range("B16").Value = 0   # crude oil
$ws.Range("B17").Value = 0   # heavy or residual fuel oil
$ws.Range("B18").Value = 0   # municipal solid waste

# "... w CCS" rows: new 95% capture share for actual CCS-equipped plant types
$ws.Range("B19").Value = 0.95   # hard coal w CCS
$ws.Range("B20").Value = 0.95   # natural gas combined cycle w CCS
$ws.Range("B21").Value = 0.95   # biomass w CCS
$ws.Range("B22").Value = 0.95   # lignite w CCS

# These two are not CCS plant types in the real sense - reset to 0
$ws.Range("B23").Value = 0   # small modular reactor
$ws.Range("B24").Value = 0   # hydrogen

# Clear the now-unused highlight style from rows 19:24 (both columns)
$ws.Range("A19:B24").ClearFormats()

# Update the sheet's active selection to reflect the new area of interest
[void]$ws.Range("B19:B22").Select()

# Restore the originally active sheet/tab (selecting a range above makes
# its sheet active as a side effect) - put the focus back on "About".
$wsAbout = $wb.Worksheets.Item("About")
[void]$wsAbout.Select()
